$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row: C/D become spearman_r/spearman_p, E becomes n, F becomes mean_diff
$ws.Range("C1").Value = "spearman_r"
$ws.Range("D1").Value = "spearman_p"
$ws.Range("E1").Value = "n"
$ws.Range("F1").Value = "mean_diff"

# New data values per row: spearman_r, spearman_p, n, mean_diff(new, recalculated)
$ws.Range("C2").Value = [double]"0.02578971968863382"
$ws.Range("D2").Value = [double]"0.6437244148461119"
$ws.Range("E2").Value = [double]"324"
$ws.Range("F2").Value = [double]"92.25689300411524"

$ws.Range("C3").Value = [double]"-0.5366818968278821"
$ws.Range("D3").Value = [double]"1.44973140734144e-25"
$ws.Range("E3").Value = [double]"324"
$ws.Range("F3").Value = [double]"112.2081275720165"

$ws.Range("C4").Value = [double]"0.07791892983253103"
$ws.Range("D4").Value = [double]"0.1617385594216551"
$ws.Range("E4").Value = [double]"324"
$ws.Range("F4").Value = [double]"1.134897119341559"

$ws.Range("C5").Value = [double]"-0.5386697465698866"
$ws.Range("D5").Value = [double]"8.90155837032734e-26"
$ws.Range("E5").Value = [double]"324"
$ws.Range("F5").Value = [double]"21.08613168724279"

$ws.Range("C6").Value = [double]"0.1600158738541625"
$ws.Range("D6").Value = [double]"0.003879661697487937"
$ws.Range("E6").Value = [double]"324"
$ws.Range("F6").Value = [double]"107.7940118312757"

$ws.Range("C7").Value = [double]"0.1369346169458027"
$ws.Range("D7").Value = [double]"0.01362777088369083"
$ws.Range("E7").Value = [double]"324"
$ws.Range("F7").Value = [double]"127.7452463991769"

$ws.Range("C8").Value = [double]"0.08165690046153645"
$ws.Range("D8").Value = [double]"0.1424866265162245"
$ws.Range("E8").Value = [double]"324"
$ws.Range("F8").Value = [double]"12145.06078549383"

$ws.Range("C9").Value = [double]"-0.5053471280814221"
$ws.Range("D9").Value = [double]"2.099657985557351e-22"
$ws.Range("E9").Value = [double]"324"
$ws.Range("F9").Value = [double]"12165.01202006173"

# Remove the now-unused columns G and H (old n/mean_diff columns no longer needed
# since data shifted left into C:F)
$ws.Range("G1:H9").Delete() | Out-Null

$wb.Save()
